$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dadosDeAcesso")

# Add new row 9 data: ID_0011, André Automatizador, sem email, automacaoteste
$ws.Cells.Item(9, 1).Value = "ID_0011"
$ws.Cells.Item(9, 2).Value = "André Automatizador"
$ws.Cells.Item(9, 3).Value = "sem email"
$ws.Cells.Item(9, 4).Value = "automacaoteste"

# Update selection to A9 as in the target workbook
$ws.Range("A9").Select()
